$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2519.3333
$ws.Range("I53").Value = 15.25
$ws.Range("J53").Value = 4522.6
$ws.Range("K53").Value = 15.25
$ws.Range("L53").Value = 4522.6
$ws.Range("M53").Value = 621.75
$ws.Range("N53").Value = -5796.6
$ws.Range("H70").Value = 257925.5
$ws.Range("I70").Value = 999999
$ws.Range("J70").Value = 10567.667
$ws.Range("K70").Value = 2999997
$ws.Range("L70").Value = 31703.001
$ws.Range("M70").Value = -2999727
$ws.Range("N70").Value = -32243.001
$ws.Range("H73").Value = 257925.5
$ws.Range("I73").Value = 999999
$ws.Range("J73").Value = 10567.667
$ws.Range("K73").Value = 2999997
$ws.Range("L73").Value = 31703.001
$ws.Range("M73").Value = -2999061
$ws.Range("N73").Value = -33575.001
$ws.Range("H76").Value = 5015.5
$ws.Range("I76").Value = 3727.5881
$ws.Range("K76").Value = 3727.5881
$ws.Range("M76").Value = -3412.5881
$ws.Range("H79").Value = 5015.5
$ws.Range("I79").Value = 3727.5881
$ws.Range("K79").Value = 3727.5881
$ws.Range("M79").Value = -2635.5881
$ws.Range("H86").Value = 4860
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4860
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4860
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -7106
$ws.Range("H89").Value = 4860
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 4860
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 24300
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -35532
$ws.Range("H98").Value = 1096.9333
$ws.Range("I98").Value = 1020
$ws.Range("K98").Value = 1020
$ws.Range("M98").Value = 478
$ws.Range("H122").Value = 1096.9333
$ws.Range("I122").Value = 1020
$ws.Range("K122").Value = 3060
$ws.Range("M122").Value = -610
$ws.Range("H138").Value = 3147.6492
$ws.Range("I138").Value = 1530.6571
$ws.Range("J138").Value = 5720.136
$ws.Range("K138").Value = 4591.971299999999
$ws.Range("L138").Value = 17160.408
$ws.Range("M138").Value = 548.0287000000008
$ws.Range("N138").Value = -27440.408

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2547.1082
$ws.Range("I32").Value = 2119.5352
$ws.Range("K32").Value = 2119.5352
$ws.Range("M32").Value = -1832.5352
$ws.Range("H43").Value = 43267.816
$ws.Range("I43").Value = 36777.332
$ws.Range("J43").Value = 45701.75
$ws.Range("K43").Value = 36777.332
$ws.Range("L43").Value = 45701.75
$ws.Range("M43").Value = -36464.332
$ws.Range("N43").Value = -46327.75
$ws.Range("H80").Value = 93332
$ws.Range("J80").Value = 93332
$ws.Range("L80").Value = 93332
$ws.Range("N80").Value = -95328
$ws.Range("H83").Value = 93332
$ws.Range("J83").Value = 93332
$ws.Range("L83").Value = 279996
$ws.Range("N83").Value = -289980
$ws.Range("H122").Value = 166669570
$ws.Range("I122").Value = 3489.8
$ws.Range("K122").Value = 10469.4
$ws.Range("M122").Value = -8019.400000000001
$ws.Range("H124").Value = 41949.715
$ws.Range("J124").Value = 41949.715
$ws.Range("L124").Value = 41949.715
$ws.Range("N124").Value = -51769.715
$ws.Range("H125").Value = 63809.668
$ws.Range("J125").Value = 63809.668
$ws.Range("L125").Value = 63809.668
$ws.Range("N125").Value = -73649.66800000001

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4427.7
$ws.Range("I86").Value = 894.75
$ws.Range("J86").Value = 6783
$ws.Range("K86").Value = 894.75
$ws.Range("L86").Value = 6783
$ws.Range("M86").Value = 228.25
$ws.Range("N86").Value = -9029
$ws.Range("H89").Value = 4427.7
$ws.Range("I89").Value = 894.75
$ws.Range("J89").Value = 6783
$ws.Range("K89").Value = 4473.75
$ws.Range("L89").Value = 33915
$ws.Range("M89").Value = 1142.25
$ws.Range("N89").Value = -45147

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34488.85
$ws.Range("I31").Value = 1735.6111
$ws.Range("K31").Value = 1735.6111
$ws.Range("M31").Value = -1440.6111
$ws.Range("H34").Value = 34488.85
$ws.Range("I34").Value = 1735.6111
$ws.Range("K34").Value = 1735.6111
$ws.Range("M34").Value = -1533.6111
$ws.Range("H58").Value = 3446.4614
$ws.Range("I58").Value = 1476.2222
$ws.Range("K58").Value = 1476.2222
$ws.Range("M58").Value = -1273.2222
$ws.Range("H112").Value = 80953.5
$ws.Range("J112").Value = 80953.5
$ws.Range("L112").Value = 80953.5
$ws.Range("N112").Value = -83907.5
$ws.Range("H136").Value = 3446.4614
$ws.Range("I136").Value = 1476.2222
$ws.Range("K136").Value = 4428.6666
$ws.Range("M136").Value = -1878.6666

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 10909.5
$ws.Range("J32").Value = 12811.4
$ws.Range("L32").Value = 38434.2
$ws.Range("N32").Value = -39000.2
$ws.Range("H74").Value = 10632.5
$ws.Range("I74").Value = 2250
$ws.Range("J74").Value = 19015
$ws.Range("K74").Value = 6750
$ws.Range("L74").Value = 57045
$ws.Range("M74").Value = -5689
$ws.Range("N74").Value = -59167
$ws.Range("H77").Value = 10632.5
$ws.Range("I77").Value = 2250
$ws.Range("J77").Value = 19015
$ws.Range("K77").Value = 20250
$ws.Range("L77").Value = 171135
$ws.Range("M77").Value = -14946
$ws.Range("N77").Value = -181743
$ws.Range("H81").Value = 3819.9092
$ws.Range("I81").Value = 1870.3334
$ws.Range("J81").Value = 5169.615
$ws.Range("K81").Value = 5611.0002
$ws.Range("L81").Value = 15508.845
$ws.Range("M81").Value = -4488.0002
$ws.Range("N81").Value = -17754.845
$ws.Range("H84").Value = 3819.9092
$ws.Range("I84").Value = 1870.3334
$ws.Range("J84").Value = 5169.615
$ws.Range("K84").Value = 16833.0006
$ws.Range("L84").Value = 46526.535
$ws.Range("M84").Value = -11217.0006
$ws.Range("N84").Value = -57758.535
$ws.Range("H87").Value = 16858.572
$ws.Range("I87").Value = 16858.572
$ws.Range("K87").Value = 50575.716
$ws.Range("M87").Value = -49327.716
$ws.Range("H88").Value = 15316.75
$ws.Range("J88").Value = 13755.667
$ws.Range("L88").Value = 41267.001
$ws.Range("N88").Value = -42123.001
$ws.Range("H90").Value = 16858.572
$ws.Range("I90").Value = 16858.572
$ws.Range("K90").Value = 151727.148
$ws.Range("M90").Value = -145487.148
$ws.Range("H91").Value = 15316.75
$ws.Range("J91").Value = 13755.667
$ws.Range("L91").Value = 41267.001
$ws.Range("N91").Value = -44231.001
$ws.Range("H97").Value = 376.54544
$ws.Range("I97").Value = 237.375
$ws.Range("J97").Value = 747.6667
$ws.Range("K97").Value = 712.125
$ws.Range("L97").Value = 2243.0001
$ws.Range("M97").Value = -216.125
$ws.Range("N97").Value = -3235.0001
$ws.Range("H121").Value = 1210
$ws.Range("J121").Value = 875.8570999999999
$ws.Range("L121").Value = 2627.5713
$ws.Range("N121").Value = -5247.5713

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 22495
$ws.Range("J44").Value = 24990
$ws.Range("L44").Value = 24990
$ws.Range("N44").Value = -26182
$ws.Range("H47").Value = 616499.5
$ws.Range("J47").Value = 616499.5
$ws.Range("L47").Value = 616499.5
$ws.Range("N47").Value = -617635.5
$ws.Range("H58").Value = 19995
$ws.Range("J58").Value = 19995
$ws.Range("L58").Value = 19995
$ws.Range("N58").Value = -20549
$ws.Range("H80").Value = 7409.3335
$ws.Range("I80").Value = 966.6667
$ws.Range("J80").Value = 9556.888999999999
$ws.Range("K80").Value = 966.6667
$ws.Range("L80").Value = 9556.888999999999
$ws.Range("M80").Value = 31.33330000000001
$ws.Range("N80").Value = -11552.889
$ws.Range("H83").Value = 7409.3335
$ws.Range("I83").Value = 966.6667
$ws.Range("J83").Value = 9556.888999999999
$ws.Range("K83").Value = 4833.3335
$ws.Range("L83").Value = 47784.44499999999
$ws.Range("M83").Value = 158.6665000000003
$ws.Range("N83").Value = -57768.44499999999
$ws.Range("H122").Value = 6942.4
$ws.Range("I122").Value = 7444.9375
$ws.Range("J122").Value = 4932.25
$ws.Range("K122").Value = 22334.8125
$ws.Range("L122").Value = 14796.75
$ws.Range("M122").Value = -19884.8125
$ws.Range("N122").Value = -19696.75
$ws.Range("H130").Value = 69549.5
$ws.Range("J130").Value = 69549.5
$ws.Range("L130").Value = 69549.5
$ws.Range("N130").Value = -79589.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8170.4
$ws.Range("I40").Value = 7580.846
$ws.Range("K40").Value = 7580.846
$ws.Range("M40").Value = -7444.846
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
$ws.Range("H130").Value = 68560.664
$ws.Range("J130").Value = 68560.664
$ws.Range("L130").Value = 68560.664
$ws.Range("N130").Value = -78600.664
$ws.Range("H139").Value = 60834.77
$ws.Range("I139").Value = 40824.75
$ws.Range("J139").Value = 69728.11
$ws.Range("K139").Value = 40824.75
$ws.Range("L139").Value = 69728.11
$ws.Range("M139").Value = -35684.75
$ws.Range("N139").Value = -80008.11

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 29714.285
$ws.Range("I55").Value = 23986.75
$ws.Range("K55").Value = 23986.75
$ws.Range("M55").Value = -23709.75
$ws.Range("H126").Value = 2302.2122
$ws.Range("I126").Value = 1864.3914
$ws.Range("J126").Value = 3309.2
$ws.Range("K126").Value = 5593.174199999999
$ws.Range("L126").Value = 9927.599999999999
$ws.Range("M126").Value = -3123.174199999999
$ws.Range("N126").Value = -14867.6
$ws.Range("H132").Value = 7364.4688
$ws.Range("I132").Value = 3412.2856
$ws.Range("J132").Value = 14909.546
$ws.Range("K132").Value = 10236.8568
$ws.Range("L132").Value = 44728.638
$ws.Range("M132").Value = -7706.856800000001
$ws.Range("N132").Value = -49788.638
$ws.Range("H136").Value = 3238.5881
$ws.Range("I136").Value = 2253.1875
$ws.Range("J136").Value = 19005
$ws.Range("K136").Value = 6759.5625
$ws.Range("L136").Value = 57015
$ws.Range("M136").Value = -4209.5625
$ws.Range("N136").Value = -62115
